# Applies the cryptos.xlsx price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cell with the default (un-styled) format, used to strip the
# text-quote-prefix style that Excel applies when a numeric-looking
# string is forced to Text below - keeps styling identical to before.
$defaultStyle = $ws.Range("B2").Style

function Set-TextCell($cell, $value) {
    # Leading apostrophe forces Excel to store the value as text even
    # when it looks like a number (e.g. "537.56"); we then reapply the
    # default style so the quote-prefix formatting does not linger.
    $rng = $ws.Range($cell)
    $rng.Value = "'" + $value
    $rng.Style = $defaultStyle
}

# Row 2
$ws.Range("D2").Value = '59.246.64'
$ws.Range("E2").Value = '  +0.37%  '

# Row 3
$ws.Range("D3").Value = '2.524.04'
$ws.Range("E3").Value = '  +0.34%  '

# Row 4
$ws.Range("E4").Value = '  -0.01%  '

# Row 5
Set-TextCell "D5" '537.56'
$ws.Range("E5").Value = '  +1.37%  '

# Row 6
Set-TextCell "D6" '138.43'
$ws.Range("E6").Value = '  -0.28%  '

# Row 7
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("E8").Value = '  +0.54%  '

# Row 9
$ws.Range("D9").Value = '2.523.03'
$ws.Range("E9").Value = '  +0.11%  '

# Row 10
$ws.Range("E10").Value = '  +1.05%  '

# Row 11
Set-TextCell "D11" '0.160'
$ws.Range("E11").Value = '  -0.19%  '

# Row 12
$ws.Range("E12").Value = '  -1.45%  '

# Row 13
$ws.Range("E13").Value = '  -2.44%  '

# Row 14
$ws.Range("D14").Value = '2.974.97'
$ws.Range("E14").Value = '  +0.50%  '

# Row 15
Set-TextCell "D15" '23.21'
$ws.Range("E15").Value = '  +0.90%  '

# Row 16
$ws.Range("D16").Value = '59.143.25'
$ws.Range("E16").Value = '  +0.28%  '

# Row 17
$ws.Range("E17").Value = '  -0.04%  '

# Row 18
$ws.Range("D18").Value = '2.520.85'
$ws.Range("E18").Value = '  +0.50%  '

# Row 19
$ws.Range("E19").Value = '  +0.90%  '

# Row 20
Set-TextCell "D20" '4.30'
$ws.Range("E20").Value = '  +1.11%  '

# Row 21
Set-TextCell "D21" '326.25'
$ws.Range("E21").Value = '  +1.20%  '

# Row 23
Set-TextCell "D23" '5.90'
$ws.Range("E23").Value = '  +1.43%  '

# Row 24
Set-TextCell "D24" '65.74'
$ws.Range("E24").Value = '  +5.60%  '

# Row 25
Set-TextCell "D25" '0.424'
$ws.Range("E25").Value = '  -0.01%  '

# Row 27
$ws.Range("E27").Value = '  +0.33%  '

# Row 28
Set-TextCell "D28" '7.66'
$ws.Range("E28").Value = '  -1.67%  '

# Row 29
Set-TextCell "D29" '6.77'
$ws.Range("E29").Value = '  -0.40%  '

# Row 30
$ws.Range("D30").Value = '0.0₃0778'
$ws.Range("E30").Value = '  +1.05%  '

# Row 31
$ws.Range("E31").Value = '  +0.22%  '

# Row 32
Set-TextCell "D32" '169.73'
$ws.Range("E32").Value = '  +4.66%  '

# Row 33
$ws.Range("E33").Value = '  +6.45%  '

# Row 34
Set-TextCell "D34" '0.999'
$ws.Range("E34").Value = '  +0.00%  '

# Row 35
$ws.Range("E35").Value = '  +2.68%  '

# Row 36
Set-TextCell "D36" '18.55'
$ws.Range("E36").Value = '  +0.53%  '

# Row 37
Set-TextCell "D37" '4.12'
$ws.Range("E37").Value = '  -2.06%  '

# Row 38
$ws.Range("E38").Value = '  -0.27%  '

# Row 39
Set-TextCell "D39" '36.69'
$ws.Range("E39").Value = '  -0.87%  '

# Row 40
Set-TextCell "D40" '0.826'
$ws.Range("E40").Value = '  +3.05%  '

# Row 41
$ws.Range("E41").Value = '  +0.33%  '

# Row 42
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell "D42" '5.27'
$ws.Range("E42").Value = '  +1.60%  '

# Row 43
$ws.Range("B43").Value = 'Bittensor'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextCell "D43" '284.28'
$ws.Range("E43").Value = '  +1.69%  '

# Row 44
$ws.Range("E44").Value = '  +0.03%  '

# Row 45
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextCell "D45" '0.606'
$ws.Range("E45").Value = '  +1.78%  '

# Row 46
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextCell "D46" '130.55'
$ws.Range("E46").Value = '  +7.18%  '

# Row 47
$ws.Range("E47").Value = '  +0.28%  '

# Row 48
$ws.Range("E48").Value = '  +0.17%  '

# Row 49
Set-TextCell "D49" '0.0511'
$ws.Range("E49").Value = '  +0.33%  '

# Row 50
$ws.Range("E50").Value = '  +0.12%  '

# Row 51
Set-TextCell "D51" '17.52'
$ws.Range("E51").Value = '  -0.27%  '
